# Update market-board derived columns (H:N) across the Leve profit sheets.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# ALC (sheet1): clear stale market data for rows 125-141
# ---------------------------------------------------------------------------
$wsALC = $wb.Worksheets.Item("ALC")
$wsALC.Range("H125:N141").ClearContents()

# ---------------------------------------------------------------------------
# ARM (sheet2): refresh row 32, clear stale market data for rows 121-141
# ---------------------------------------------------------------------------
$wsARM = $wb.Worksheets.Item("ARM")
$wsARM.Range("H32").Value2 = 1145040.2
$wsARM.Range("I32").Value2 = 1193902
$wsARM.Range("K32").Value2 = 1193902
$wsARM.Range("M32").Value2 = -1193615
$wsARM.Range("H121:N141").ClearContents()

# ---------------------------------------------------------------------------
# BSM (sheet3): refresh rows 86 and 89
# ---------------------------------------------------------------------------
$wsBSM = $wb.Worksheets.Item("BSM")
$wsBSM.Range("H86").Value2 = 6084.7856
$wsBSM.Range("I86").Value2 = 7125.5454
$wsBSM.Range("J86").Value2 = 2268.6667
$wsBSM.Range("K86").Value2 = 7125.5454
$wsBSM.Range("L86").Value2 = 2268.6667
$wsBSM.Range("M86").Value2 = -6002.5454
$wsBSM.Range("N86").Value2 = -4514.6667

$wsBSM.Range("H89").Value2 = 6084.7856
$wsBSM.Range("I89").Value2 = 7125.5454
$wsBSM.Range("J89").Value2 = 2268.6667
$wsBSM.Range("K89").Value2 = 35627.727
$wsBSM.Range("L89").Value2 = 11343.3335
$wsBSM.Range("M89").Value2 = -30011.727
$wsBSM.Range("N89").Value2 = -22575.3335

# ---------------------------------------------------------------------------
# CRP (sheet4): refresh rows 99, 126 and 132
# ---------------------------------------------------------------------------
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCRP.Range("H99").Value2 = 33776.098
$wsCRP.Range("I99").Value2 = 46816.316
$wsCRP.Range("J99").Value2 = 1900
$wsCRP.Range("K99").Value2 = 46816.316
$wsCRP.Range("L99").Value2 = 1900
$wsCRP.Range("M99").Value2 = -45318.316
$wsCRP.Range("N99").Value2 = -4896

$wsCRP.Range("H126").Value2 = 33776.098
$wsCRP.Range("I126").Value2 = 46816.316
$wsCRP.Range("J126").Value2 = 1900
$wsCRP.Range("K126").Value2 = 140448.948
$wsCRP.Range("L126").Value2 = 5700
$wsCRP.Range("M126").Value2 = -137978.948
$wsCRP.Range("N126").Value2 = -10640

$wsCRP.Range("H132").Value2 = 11365395
$wsCRP.Range("I132").Value2 = 20834534
$wsCRP.Range("J132").Value2 = 2428.85
$wsCRP.Range("K132").Value2 = 62503602
$wsCRP.Range("L132").Value2 = 7286.549999999999
$wsCRP.Range("M132").Value2 = -62501072
$wsCRP.Range("N132").Value2 = -12346.55

# ---------------------------------------------------------------------------
# CUL (sheet5): refresh rows 68 and 71, clear stale market data for 120-141
# ---------------------------------------------------------------------------
$wsCUL = $wb.Worksheets.Item("CUL")
$wsCUL.Range("H68").Value2 = 1327.317
$wsCUL.Range("J68").Value2 = 1606.5186
$wsCUL.Range("L68").Value2 = 4819.5558
$wsCUL.Range("N68").Value2 = -6441.5558

$wsCUL.Range("H71").Value2 = 1327.317
$wsCUL.Range("J71").Value2 = 1606.5186
$wsCUL.Range("L71").Value2 = 14458.6674
$wsCUL.Range("N71").Value2 = -22570.6674

$wsCUL.Range("H120:N141").ClearContents()

# ---------------------------------------------------------------------------
# LTW (sheet7): refresh rows 68, 71, 88 and 91
# ---------------------------------------------------------------------------
$wsLTW = $wb.Worksheets.Item("LTW")
$wsLTW.Range("H68").Value2 = 2223.5557
$wsLTW.Range("I68").Value2 = 2113.5789
$wsLTW.Range("J68").Value2 = 2484.75
$wsLTW.Range("K68").Value2 = 2113.5789
$wsLTW.Range("L68").Value2 = 2484.75
$wsLTW.Range("M68").Value2 = -1364.5789
$wsLTW.Range("N68").Value2 = -3982.75

$wsLTW.Range("H71").Value2 = 2223.5557
$wsLTW.Range("I71").Value2 = 2113.5789
$wsLTW.Range("J71").Value2 = 2484.75
$wsLTW.Range("K71").Value2 = 10567.8945
$wsLTW.Range("L71").Value2 = 12423.75
$wsLTW.Range("M71").Value2 = -6823.8945
$wsLTW.Range("N71").Value2 = -19911.75

$wsLTW.Range("H88").Value2 = 35500
$wsLTW.Range("J88").Value2 = 35500
$wsLTW.Range("L88").Value2 = 35500
$wsLTW.Range("N88").Value2 = -36356

$wsLTW.Range("H91").Value2 = 35500
$wsLTW.Range("J91").Value2 = 35500
$wsLTW.Range("L91").Value2 = 35500
$wsLTW.Range("N91").Value2 = -38464

# ---------------------------------------------------------------------------
# WVR (sheet8): refresh row 136
# ---------------------------------------------------------------------------
$wsWVR = $wb.Worksheets.Item("WVR")
$wsWVR.Range("H136").Value2 = 22004162
$wsWVR.Range("I136").Value2 = 28890316
$wsWVR.Range("J136").Value2 = 6940700.5
$wsWVR.Range("K136").Value2 = 86670948
$wsWVR.Range("L136").Value2 = 20822101.5
$wsWVR.Range("M136").Value2 = -86668398
$wsWVR.Range("N136").Value2 = -20827201.5
